$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 439, shifting existing rows 439:540 down to 440:541
$ws.Rows.Item(439).Insert()

# Populate the newly inserted row 439 with the new record
$ws.Cells.Item(439, 1).Value = 10
$ws.Cells.Item(439, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(439, 3).Value = "La Araucanía"
$ws.Cells.Item(439, 4).Value = 44508
$ws.Cells.Item(439, 5).Value = 9
$ws.Cells.Item(439, 6).Value = "Fruta"
$ws.Cells.Item(439, 7).Value = 100109
$ws.Cells.Item(439, 8).Value = "Uva"
$ws.Cells.Item(439, 9).Value = 100109001
$ws.Cells.Item(439, 10).Value = "Uva"
$ws.Cells.Item(439, 11).Value = "Superior Seedless"
$ws.Cells.Item(439, 12).Value = "Primera"
$ws.Cells.Item(439, 13).Value = 280
$ws.Cells.Item(439, 14).Value = 33000
$ws.Cells.Item(439, 15).Value = 33000
$ws.Cells.Item(439, 16).Value = 33000
$ws.Cells.Item(439, 17).Value = "$/bandeja 8 kilos"
$ws.Cells.Item(439, 18).Value = "EE.UU."
$ws.Cells.Item(439, 19).Value = 4125
$ws.Cells.Item(439, 20).Value = 8
